# Reorder the worker/period detail rows (B16:G25) into the new order.
# The set of (worker, period, value, salary) tuples is unchanged; only the
# row order changes, grouping DAVID RICARDO POLO PALENCIA's four periods
# first (descending 2304->2301), then the existing CRISTIAN/LEONARDO rows,
# then DIEGO ANDRES HERNANDEZ LUNA's four periods (descending 2304->2301).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$endRow = 25

# Snapshot current B:G values for each row in the block.
$rows = @{}
for ($r = $startRow; $r -le $endRow; $r++) {
    $rows[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2
    )
}

# New row order expressed as the source row (within the old snapshot) that
# should now occupy each destination row 16..25.
$newOrder = @(24, 23, 20, 18, 16, 17, 25, 22, 21, 19)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $destRow = $startRow + $i
    $srcRow = $newOrder[$i]
    $vals = $rows[$srcRow]

    $ws.Cells.Item($destRow, 2).Value2 = $vals[0]
    $ws.Cells.Item($destRow, 3).Value2 = $vals[1]
    $ws.Cells.Item($destRow, 4).Value2 = $vals[2]
    $ws.Cells.Item($destRow, 5).Value2 = $vals[3]
    $ws.Cells.Item($destRow, 6).Value2 = $vals[4]
    $ws.Cells.Item($destRow, 7).Value2 = $vals[5]
}
